# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Julio de 2020 a las 05:03"

# --- Reorder: Honduras now ranks above Guatemala (rows 55/56 swap names) ---
$ws.Range("A55").Value = "Honduras"
$ws.Range("A56").Value = "Guatemala"

# --- Reorder: Islas Malvinas now ranks above Groenlandia (rows 209/210 swap names) ---
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Brasil (row 5): updated Casos activos / Recuperados ---
$ws.Range("D5").Value = 1117922
$ws.Range("E5").Value = 489865

# --- Bolivia (row 43): updated stats ---
$ws.Range("B43").Value = 41545
$ws.Range("C43").Value = 1036
$ws.Range("D43").Value = 12398
$ws.Range("E43").Value = 27617
$ws.Range("G43").Value = 54
$ws.Range("H43").Value = 1530

# --- Row 55 (now Honduras): updated stats ---
$ws.Range("B55").Value = 25428
$ws.Range("C55").Value = 763
$ws.Range("D55").Value = 2637
$ws.Range("E55").Value = 22114
$ws.Range("G55").Value = 21
$ws.Range("H55").Value = 677

# --- Row 56 (now Guatemala): updated stats ---
$ws.Range("B56").Value = 24787
$ws.Range("D56").Value = 3575
$ws.Range("E56").Value = 20208
$ws.Range("H56").Value = 1004

# --- Nueva Zelanda (row 125): updated stats ---
$ws.Range("B125").Value = 1537
$ws.Range("C125").Value = 1
$ws.Range("E125").Value = 23
